$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.252593755722046
$ws.Range("B1").Value = 2.634604454040527
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.697762370109558
$ws.Range("E1").Value = 1.134664297103882
